# Add thêm nhân sự Nguyễn Hữu Quang
# Updates the "Lương" (Salary) sheet: one additional work day at LONG XUYÊN,
# which increases base/total salary proportionally.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Danh mục lương (count of salary line items): 20 -> 21
$ws.Range("B1").Value = 21

# Tổng công tại LONG XUYÊN (total work days at LONG XUYEN): 28 -> 29
$ws.Range("B12").Value = 29

# Lương cơ bản tại LONG XUYÊN (base salary at LONG XUYEN)
$ws.Range("B13").Value = 5178571.428571429

# Tổng lương tại LONG XUYÊN (total salary at LONG XUYEN)
$ws.Range("B32").Value = 5178571.428571429

# Tổng lương tại HỆ THỐNG (total salary system-wide)
$ws.Range("B34").Value = 6118571.428571429
